# Update volume/price/origin data and correct dates for the weekly
# fruit/vegetable price refresh (Maracuyá - Vega Modelo de Temuco).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44424
$ws.Range("M2").Value = 15
$ws.Range("R2").Value = 'Región de Arica y Parinacota'

$ws.Range("D3").Value = 44294
$ws.Range("M3").Value = 15

$ws.Range("D4").Value = 44364
$ws.Range("M4").Value = 90
$ws.Range("N4").Value = 1700
$ws.Range("O4").Value = 1700
$ws.Range("P4").Value = 1700
$ws.Range("Q4").Value = '$/kilo'
$ws.Range("S4").Value = 1700
$ws.Range("T4").Value = 1

$ws.Range("D5").Value = 44379
$ws.Range("M5").Value = 10
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("S5").Value = 1667

$ws.Range("D6").Value = 44418
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 35000
$ws.Range("O6").Value = 35000
$ws.Range("P6").Value = 35000
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("S6").Value = 1944
$ws.Range("T6").Value = 18

$ws.Range("D7").Value = 44377
$ws.Range("N7").Value = 40000
$ws.Range("O7").Value = 40000
$ws.Range("P7").Value = 40000
$ws.Range("S7").Value = 2222

$ws.Range("D8").Value = 44369
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = 35000
$ws.Range("O8").Value = 35000
$ws.Range("P8").Value = 35000
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("R8").Value = 'Perú'
$ws.Range("S8").Value = 1944
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44433
$ws.Range("M9").Value = 15

$ws.Range("D10").Value = 44438
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 35000
$ws.Range("O10").Value = 35000
$ws.Range("P10").Value = 35000
$ws.Range("S10").Value = 1944

$ws.Range("D11").Value = 44392
$ws.Range("M11").Value = 20
$ws.Range("R11").Value = 'Región de Arica y Parinacota'

$ws.Range("D12").Value = 44434
$ws.Range("M12").Value = 40

$ws.Range("D13").Value = 44449
$ws.Range("M13").Value = 20
$ws.Range("N13").Value = 38000
$ws.Range("O13").Value = 38000
$ws.Range("P13").Value = 38000
$ws.Range("S13").Value = 2111

$ws.Range("D14").Value = 44405
$ws.Range("M14").Value = 10
$ws.Range("N14").Value = 35000
$ws.Range("O14").Value = 35000
$ws.Range("P14").Value = 35000
$ws.Range("S14").Value = 1944

$ws.Range("D15").Value = 44435
$ws.Range("R15").Value = 'Perú'

$ws.Range("D16").Value = 44435
$ws.Range("M16").Value = 105
$ws.Range("R16").Value = 'Región de Arica y Parinacota'

$ws.Range("D17").Value = 44442
$ws.Range("R17").Value = 'Perú'

$ws.Range("D18").Value = 44357
$ws.Range("M18").Value = 10
$ws.Range("N18").Value = 38000
$ws.Range("O18").Value = 38000
$ws.Range("P18").Value = 38000
$ws.Range("R18").Value = 'Perú'
$ws.Range("S18").Value = 2111

$ws.Range("D19").Value = 44264
$ws.Range("M19").Value = 20
$ws.Range("N19").Value = 40000
$ws.Range("O19").Value = 40000
$ws.Range("P19").Value = 40000
$ws.Range("R19").Value = 'Región de Arica y Parinacota'
$ws.Range("S19").Value = 2222

$ws.Range("D20").Value = 44363
$ws.Range("M20").Value = 144
$ws.Range("N20").Value = 1700
$ws.Range("O20").Value = 1700
$ws.Range("P20").Value = 1700
$ws.Range("Q20").Value = '$/kilo'
$ws.Range("S20").Value = 1700
$ws.Range("T20").Value = 1

$ws.Range("D21").Value = 44431
$ws.Range("M21").Value = 30
$ws.Range("R21").Value = 'Región de Arica y Parinacota'

$ws.Range("D22").Value = 44448
$ws.Range("M22").Value = 50

$ws.Range("D23").Value = 44279
$ws.Range("M23").Value = 30
$ws.Range("O23").Value = 36000
$ws.Range("P23").Value = 35667
$ws.Range("S23").Value = 1982

$ws.Range("D24").Value = 44432
$ws.Range("M24").Value = 10
$ws.Range("R24").Value = 'Perú'
